$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 19.3660239309531
$ws.Range("G2").Value = 56.4246557170697
$ws.Range("H2").Value = 1.24002776769199
$ws.Range("E3").Value = 32.3115930649242
$ws.Range("G3").Value = 71.1491096173724
$ws.Range("H3").Value = 2.75646234724735
$ws.Range("E4").Value = 12.1487339167937
$ws.Range("G4").Value = 46.2554251789948
$ws.Range("H4").Value = 1.51551188462199
$ws.Range("E7").Value = 2.55244410098195
$ws.Range("G7").Value = 50.9253433843927
$ws.Range("H7").Value = 0.208498319480921
$ws.Range("E8").Value = 1.08160773776441
$ws.Range("G8").Value = 159.834098987101
$ws.Range("H8").Value = 0.0626052595076097
$ws.Range("E9").Value = 11.4508287546717
$ws.Range("G9").Value = 88.4059258706094
$ws.Range("H9").Value = 0.75495402656921
$ws.Range("E10").Value = 10.6882507036427
$ws.Range("G10").Value = 57.3765946510106
$ws.Range("H10").Value = 0.668878174634346
$ws.Range("E11").Value = 13.4822442765145
$ws.Range("G11").Value = 80.822234425004
$ws.Range("H11").Value = 0.721154780964922
$ws.Range("E12").Value = 2.09171913182225
$ws.Range("G12").Value = 63.6180566319126
$ws.Range("H12").Value = 0.226875357049624
$ws.Range("E13").Value = 7.14184929086677
$ws.Range("G13").Value = 60.4552559663563
$ws.Range("H13").Value = 0.397586236678635
$ws.Range("E14").Value = 7.96029094711111
$ws.Range("G14").Value = 111.832173558297
$ws.Range("H14").Value = 0.0978407010214374
$ws.Range("E15").Value = 4.95424749782809
$ws.Range("G15").Value = 91.2398287305991
$ws.Range("H15").Value = 0.112526447360513
$ws.Range("E16").Value = 9.90883240205773
$ws.Range("G16").Value = 79.1350533711568
$ws.Range("H16").Value = 0.613594647819945
$ws.Range("E17").Value = 6.89787302636999
$ws.Range("G17").Value = 81.8326142577399
$ws.Range("H17").Value = 0.4572157936862
$ws.Range("E18").Value = 1.9752706239993
$ws.Range("G18").Value = 108.712219396557
$ws.Range("H18").Value = 0.0403439357801364
$ws.Range("E19").Value = 1.94359723892004
$ws.Range("G19").Value = 159.98439860285
$ws.Range("H19").Value = 0.0803180753654274
$ws.Range("E20").Value = 1.10171660765841
$ws.Range("G20").Value = 56.4179707284158
$ws.Range("H20").Value = 0.0952979149218019
$ws.Range("E22").Value = 8.05543022049837
$ws.Range("G22").Value = 38.639548748505
$ws.Range("H22").Value = 1.15134474568339
$ws.Range("E23").Value = 4.06009105428833
$ws.Range("G23").Value = 62.1929400577247
$ws.Range("H23").Value = 0.312925706503849
$ws.Range("E24").Value = 10.3087541506696
$ws.Range("G24").Value = 45.4471516020832
$ws.Range("H24").Value = 1.18275181904848
$ws.Range("E25").Value = 16.2740549278997
$ws.Range("G25").Value = 80.2616579862731
$ws.Range("H25").Value = 0.89987396222228
$ws.Range("E26").Value = 20.5481514883504
$ws.Range("G26").Value = 107.17562240175
$ws.Range("H26").Value = 0.677985800420298
$ws.Range("E27").Value = 11.4710898977279
$ws.Range("G27").Value = 145.80986674202
$ws.Range("H27").Value = 0.568073044269495
$ws.Range("E28").Value = 15.9950720733092
$ws.Range("G28").Value = 90.9624477341969
$ws.Range("H28").Value = 0.527852114685842
$ws.Range("E29").Value = 1.31130411886539
$ws.Range("G29").Value = 83.6509358975027
$ws.Range("H29").Value = 0.0356980576086593
$ws.Range("E30").Value = 2.50363275666084
$ws.Range("G30").Value = 199.790698543052
$ws.Range("H30").Value = 0.0427241491771618
$ws.Range("E31").Value = 8.05880367762796
$ws.Range("G31").Value = 48.3382214399045
$ws.Range("H31").Value = 0.679937559589286
$ws.Range("E33").Value = 1.00511483841174
$ws.Range("G33").Value = 80.5287841121349
$ws.Range("H33").Value = 0.118214391958584
$ws.Range("E34").Value = 7.77823614913313
$ws.Range("G34").Value = 41.6501167170211
$ws.Range("H34").Value = 1.36399706330909
$ws.Range("E35").Value = 12.7264646138699
$ws.Range("G35").Value = 100.57090441306
$ws.Range("H35").Value = 0.58735516338987
$ws.Range("E36").Value = 4.85739375472713
$ws.Range("G36").Value = 69.5327489128583
$ws.Range("H36").Value = 0.327552966036586
